$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.136.28"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.382.61"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Formula = "=`"303.54`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Formula = "=`"97.81`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").Formula = "=`"34.23`""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Formula = "=`"0.122`""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Formula = "=`"18.53`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "2.757.16"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "2.386.33"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Formula = "=`"0.811`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").Value = "43.133.17"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Formula = "=`"12.32`""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Formula = "=`"68.38`""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Formula = "=`"236.19`""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Formula = "=`"24.77`""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").Formula = "=`"2.37`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Formula = "=`"31.69`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Formula = "=`"0.0735`""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("D34").Formula = "=`"17.26`""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  +7.01%  "
$ws.Range("D36").Formula = "=`"4.38`""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Formula = "=`"0.102`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Formula = "=`"2.30`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  +4.38%  "
$ws.Range("D40").Formula = "=`"22.45`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +11.24%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Formula = "=`"107.52`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -35.10%  "
$ws.Range("D43").Value = "1.952.15"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Formula = "=`"9.26`""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -10.92%  "
$ws.Range("D48").Value = "2.618.64"
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Formula = "=`"1.51`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Formula = "=`"72.22`""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +1.54%  "

$excel.CutCopyMode = 0
